# Natmi following Dr Hou advice
# Extend Selplg-Selp sheet1 data from 3 rows (ECs/FAPs/sCs x 1 target) to a full 3x3 cluster matrix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Selplg"
$ws.Cells.Item(2,3).Value = "Selp"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 88.31721
$ws.Cells.Item(2,8).Value = 264.95163
$ws.Cells.Item(2,9).Value = 0.9709320988649861
$ws.Cells.Item(2,10).Value = 0.970932098864986
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 26.450162
$ws.Cells.Item(2,14).Value = 79.350486
$ws.Cells.Item(2,15).Value = 0.9821351879331711
$ws.Cells.Item(2,16).Value = 0.9821351879331711
$ws.Cells.Item(2,17).Value = 2336.00451188802
$ws.Cells.Item(2,18).Value = 21024.04060699218
$ws.Cells.Item(2,19).Value = 0.9535865793891114
$ws.Cells.Item(2,20).Value = 0.9535865793891113

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Selplg"
$ws.Cells.Item(3,3).Value = "Selp"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 88.31721
$ws.Cells.Item(3,8).Value = 264.95163
$ws.Cells.Item(3,9).Value = 0.9709320988649861
$ws.Cells.Item(3,10).Value = 0.970932098864986
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.321934
$ws.Cells.Item(3,14).Value = 0.965802
$ws.Cells.Item(3,15).Value = 0.01195390446349922
$ws.Cells.Item(3,16).Value = 0.01195390446349922
$ws.Cells.Item(3,17).Value = 28.43231268414
$ws.Cells.Item(3,18).Value = 255.89081415726
$ws.Cells.Item(3,19).Value = 0.01160642955037682
$ws.Cells.Item(3,20).Value = 0.01160642955037683

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Selplg"
$ws.Cells.Item(4,3).Value = "Selp"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 88.31721
$ws.Cells.Item(4,8).Value = 264.95163
$ws.Cells.Item(4,9).Value = 0.9709320988649861
$ws.Cells.Item(4,10).Value = 0.970932098864986
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.1591883333333333
$ws.Cells.Item(4,14).Value = 0.477565
$ws.Cells.Item(4,15).Value = 0.005910907603329674
$ws.Cells.Item(4,16).Value = 0.005910907603329674
$ws.Cells.Item(4,17).Value = 14.05906946455
$ws.Cells.Item(4,18).Value = 126.53162518095
$ws.Cells.Item(4,19).Value = 0.005739089925497885
$ws.Cells.Item(4,20).Value = 0.005739089925497885

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Selplg"
$ws.Cells.Item(5,3).Value = "Selp"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.137824333333333
$ws.Cells.Item(5,8).Value = 6.413473
$ws.Cells.Item(5,9).Value = 0.02350257970069449
$ws.Cells.Item(5,10).Value = 0.02350257970069449
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 26.450162
$ws.Cells.Item(5,14).Value = 79.350486
$ws.Cells.Item(5,15).Value = 0.9821351879331711
$ws.Cells.Item(5,16).Value = 0.9821351879331711
$ws.Cells.Item(5,17).Value = 56.54579994420867
$ws.Cells.Item(5,18).Value = 508.912199497878
$ws.Cells.Item(5,19).Value = 0.02308271053125592
$ws.Cells.Item(5,20).Value = 0.02308271053125592

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Selplg"
$ws.Cells.Item(6,3).Value = "Selp"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.137824333333333
$ws.Cells.Item(6,8).Value = 6.413473
$ws.Cells.Item(6,9).Value = 0.02350257970069449
$ws.Cells.Item(6,10).Value = 0.02350257970069449
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.321934
$ws.Cells.Item(6,14).Value = 0.965802
$ws.Cells.Item(6,15).Value = 0.01195390446349922
$ws.Cells.Item(6,16).Value = 0.01195390446349922
$ws.Cells.Item(6,17).Value = 0.6882383389273332
$ws.Cells.Item(6,18).Value = 6.194145050346
$ws.Cells.Item(6,19).Value = 0.000280947592387878
$ws.Cells.Item(6,20).Value = 0.0002809475923878781

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Selplg"
$ws.Cells.Item(7,3).Value = "Selp"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.137824333333333
$ws.Cells.Item(7,8).Value = 6.413473
$ws.Cells.Item(7,9).Value = 0.02350257970069449
$ws.Cells.Item(7,10).Value = 0.02350257970069449
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.1591883333333333
$ws.Cells.Item(7,14).Value = 0.477565
$ws.Cells.Item(7,15).Value = 0.005910907603329674
$ws.Cells.Item(7,16).Value = 0.005910907603329674
$ws.Cells.Item(7,17).Value = 0.3403166925827777
$ws.Cells.Item(7,18).Value = 3.062850233245
$ws.Cells.Item(7,19).Value = 0.0001389215770506967
$ws.Cells.Item(7,20).Value = 0.0001389215770506967

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Selplg"
$ws.Cells.Item(8,3).Value = "Selp"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.5062286666666665
$ws.Cells.Item(8,8).Value = 1.518686
$ws.Cells.Item(8,9).Value = 0.005565321434319426
$ws.Cells.Item(8,10).Value = 0.005565321434319426
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 26.450162
$ws.Cells.Item(8,14).Value = 79.350486
$ws.Cells.Item(8,15).Value = 0.9821351879331711
$ws.Cells.Item(8,16).Value = 0.9821351879331711
$ws.Cells.Item(8,17).Value = 13.38983024237733
$ws.Cells.Item(8,18).Value = 120.508472181396
$ws.Cells.Item(8,19).Value = 0.005465898012803815
$ws.Cells.Item(8,20).Value = 0.005465898012803815

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Selplg"
$ws.Cells.Item(9,3).Value = "Selp"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.5062286666666665
$ws.Cells.Item(9,8).Value = 1.518686
$ws.Cells.Item(9,9).Value = 0.005565321434319426
$ws.Cells.Item(9,10).Value = 0.005565321434319426
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.321934
$ws.Cells.Item(9,14).Value = 0.965802
$ws.Cells.Item(9,15).Value = 0.01195390446349922
$ws.Cells.Item(9,16).Value = 0.01195390446349922
$ws.Cells.Item(9,17).Value = 0.1629722195746666
$ws.Cells.Item(9,18).Value = 1.466749976172
$ws.Cells.Item(9,19).Value = 0.00006652732073451887
$ws.Cells.Item(9,20).Value = 0.00006652732073451888

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Selplg"
$ws.Cells.Item(10,3).Value = "Selp"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.5062286666666665
$ws.Cells.Item(10,8).Value = 1.518686
$ws.Cells.Item(10,9).Value = 0.005565321434319426
$ws.Cells.Item(10,10).Value = 0.005565321434319426
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.1591883333333333
$ws.Cells.Item(10,14).Value = 0.477565
$ws.Cells.Item(10,15).Value = 0.005910907603329674
$ws.Cells.Item(10,16).Value = 0.005910907603329674
$ws.Cells.Item(10,17).Value = 0.08058569773222221
$ws.Cells.Item(10,18).Value = 0.7252712795899999
$ws.Cells.Item(10,19).Value = 0.0000328961007810923
$ws.Cells.Item(10,20).Value = 0.0000328961007810923

